# Added Test Data for UK Market
# Copy the "Poland" worksheet (the last country sheet / template) to the end
# of the workbook, rename it to "UK", and update its market-specific cells.

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Duplicate the sheet, placing the copy right after "Poland" (i.e. at the end).
$poland.Copy($null, $poland)

# The newly created sheet becomes the active sheet after Copy(), grab it.
$uk = $wb.ActiveSheet
$uk.Name = "UK"

# Update the market-specific values that differ from the Poland template.
# (Write B4 before B2 so the new shared-string entries land in the same
# order as the target workbook: NGC-2741/T3350/T3368 then UK Market.)
$uk.Range("B4").Value = "NGC-2741/T3350/T3368"
$uk.Range("B2").Value = "UK Market"

# Match the selection recorded for the new tab in the saved workbook.
$uk.Range("K10").Select() | Out-Null
